$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overall": update row 2 stats (B2:K2)
# ----------------------------------------------------------------------
$overall = $wb.Worksheets.Item("Overall")

$overall.Range("B2").Value = 146
$overall.Range("C2").Value = 6
$overall.Range("D2").Value = 1.9521901716272532
$overall.Range("E2").Value = 0.26785244902410543
$overall.Range("F2").Value = 2.6585253456221216
$overall.Range("G2").Value = 287
$overall.Range("H2").Value = 73
$overall.Range("I2").Value = 360
$overall.Range("J2").Value = 1138
$overall.Range("K2").Value = 200

# ----------------------------------------------------------------------
# Sheet "Zones": update rows 2-14 (B:F)
# ----------------------------------------------------------------------
$zones = $wb.Worksheets.Item("Zones")

# Row 2 (Zone 1)
$zones.Range("B2").Value = 21
$zones.Range("C2").Value = 0
$zones.Range("D2").Value = 0.93275862068965509
$zones.Range("E2").Value = 0.25289855072463768
$zones.Range("F2").Value = 3.5388888888888883

# Row 3 (Zone 2)
$zones.Range("B3").Value = 10
$zones.Range("C3").Value = 1
$zones.Range("D3").Value = 2.0765151515151516
$zones.Range("E3").Value = 0.4166666666666668
$zones.Range("F3").Value = 2.4453703703703709

# Row 4 (Zone 3)
$zones.Range("B4").Value = 4
$zones.Range("C4").Value = 0
$zones.Range("D4").Value = 1.7222222222222221
$zones.Range("E4").Value = 0.23333333333333328
$zones.Range("F4").Value = 2.2636363636363637

# Row 5 (Zone 4)
$zones.Range("B5").Value = 11
$zones.Range("C5").Value = 0
$zones.Range("D5").Value = 2.2646666666666664
$zones.Range("E5").Value = 0.27333333333333359
$zones.Range("F5").Value = 2.7624999999999997

# Row 6 (Zone 5)
$zones.Range("B6").Value = 20
$zones.Range("C6").Value = 0
$zones.Range("D6").Value = 1.0075757575757578
$zones.Range("E6").Value = 0.1499999999999998
$zones.Range("F6").Value = 1.2598039215686276

# Row 7 (Zone 6)
$zones.Range("B7").Value = 15
$zones.Range("C7").Value = 0
$zones.Range("D7").Value = 1.3674242424242424
$zones.Range("E7").Value = 0.21333333333333329
$zones.Range("F7").Value = 1.7068627450980394

# Row 8 (Zone 7)
$zones.Range("B8").Value = 7
$zones.Range("C8").Value = 0
$zones.Range("D8").Value = 2.6716666666666669
$zones.Range("E8").Value = 0.26000000000000006
$zones.Range("F8").Value = 3.4755555555555562

# Row 9 (Zone 8)
$zones.Range("B9").Value = 7
$zones.Range("C9").Value = 0
$zones.Range("D9").Value = 2.4746031746031738
$zones.Range("E9").Value = 0.23666666666666672
$zones.Range("F9").Value = 3.1739583333333328

# Row 10 (Zone 9)
$zones.Range("B10").Value = 5
$zones.Range("C10").Value = 3
$zones.Range("D10").Value = 3.4172839506172838
$zones.Range("E10").Value = 0.45208333333333361
$zones.Range("F10").Value = 4.6657894736842103

# Row 11 (Zone 10)
$zones.Range("B11").Value = 8
$zones.Range("C11").Value = 0
$zones.Range("D11").Value = 3.0320000000000009
$zones.Range("E11").Value = 0.28750000000000003
$zones.Range("F11").Value = 3.5547619047619059

# Row 12 (Zone 11)
$zones.Range("B12").Value = 15
$zones.Range("C12").Value = 0
$zones.Range("D12").Value = 1.138095238095238
$zones.Range("E12").Value = 0.14666666666666678
$zones.Range("F12").Value = 1.4479166666666665

# Row 13 (Zone 12)
$zones.Range("B13").Value = 10
$zones.Range("C13").Value = 2
$zones.Range("D13").Value = 1.1291666666666662
$zones.Range("E13").Value = 0.36166666666666653
$zones.Range("F13").Value = 1.6773809523809518

# Row 14 (Zone 13)
$zones.Range("B14").Value = 13
$zones.Range("C14").Value = 0
$zones.Range("D14").Value = 1.9842857142857144
$zones.Range("E14").Value = 0.15625000000000003
$zones.Range("F14").Value = 2.5259259259259261
